# Picked out components for nutrient management
# Adds two new BOM rows (51, 52) for a water temperature sensor and a
# nutrient solution heater, with formulas, hyperlinks and formatting that
# match the existing table rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Copy formatting (column D wrap-text, column F hyperlink style,
#     columns J/K currency) from the last existing data row (50) down
#     into the two new rows, matching the look of the rest of the table.
$ws.Range("B50:K50").Copy() | Out-Null
$ws.Range("B51:K51").PasteSpecial(-4122) | Out-Null
$ws.Range("B52:K52").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Both new rows wrap their long description text, same as the row above.
$ws.Rows.Item(51).RowHeight = 30
$ws.Rows.Item(52).RowHeight = 30

# --- Row 51: Water temperature sensor -------------------------------
$ws.Range("B51").Value = "Water temperature sensor"
$ws.Range("C51").Value = "Aideepen"
$ws.Range("D51").Value = "Aideepen 5 x DS18B20 2m Cable Temperature Digital Thermal Probe Sensor Stainless Steel Probe, Accurate Reading, Measure Temperature -55°C to +125°C"
$ws.Range("E51").Value = "A70110795UK"
$ws.Range("F51").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("F51"), "https://www.amazon.co.uk/dp/B07N2YZ2NR", "", "", "Amazon") | Out-Null
# Adding the hyperlink re-styles the cell with a fresh "Hyperlink" xf;
# put back the same style already used by the rest of column F.
$ws.Range("F51").Style = $ws.Range("F50").Style
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 5
$ws.Range("I51").Formula = "=G51*H51"
$ws.Range("J51").Value = 11.99
$ws.Range("K51").Formula = "=G51*J51"

# --- Row 52: Nutrient solution heater ---------------------------------
# (Model # is entered before the product name, matching the order the
# two new shared-string table entries were originally authored in.)
$ws.Range("B52").Value = "Nutrient solution heater"
$ws.Range("C52").Value = "Hidom"
$ws.Range("E52").Value = [char]0x200e + "HT-2025-U"
$ws.Range("D52").Value = "Hidom 25w Submersible Aquarium Fish Tank Shockproof Thermostat Heater - Adjustable Temperature HT-2025"
$ws.Range("F52").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("F52"), "https://www.amazon.co.uk/dp/B0851ZHB9J", "", "", "Amazon") | Out-Null
$ws.Range("F52").Style = $ws.Range("F50").Style
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 1
$ws.Range("I52").Formula = "=G52*H52"
$ws.Range("J52").Value = 10.45
$ws.Range("K52").Formula = "=G52*J52"

# --- View state: scroll the frozen pane down to show the new rows and
#     leave the selection on the last-entered cell, like the author did.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 45
$ws.Range("H52").Select()
